# Auto-applied market price / profit snapshot refresh
# (mirrors a scheduled runner that re-pulls current Market Board prices
#  and recomputes the NQ/HQ profit columns H:N for specific Leve rows)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Cells.Item(64, 8).Value = 3851.75
$ws.Cells.Item(64, 9).Value = 3689.25
$ws.Cells.Item(64, 10).Value = 4014.25
$ws.Cells.Item(64, 11).Value = 3689.25
$ws.Cells.Item(64, 12).Value = 4014.25
$ws.Cells.Item(64, 13).Value = -3441.25
$ws.Cells.Item(64, 14).Value = -4510.25

# Row 67
$ws.Cells.Item(67, 8).Value = 3851.75
$ws.Cells.Item(67, 9).Value = 3689.25
$ws.Cells.Item(67, 10).Value = 4014.25
$ws.Cells.Item(67, 11).Value = 3689.25
$ws.Cells.Item(67, 12).Value = 4014.25
$ws.Cells.Item(67, 13).Value = -2831.25
$ws.Cells.Item(67, 14).Value = -5730.25

# Row 98
$ws.Cells.Item(98, 8).Value = 1717.6666
$ws.Cells.Item(98, 10).Value = 2611.8333
$ws.Cells.Item(98, 12).Value = 2611.8333
$ws.Cells.Item(98, 14).Value = -5607.8333

# Row 122
$ws.Cells.Item(122, 8).Value = 1717.6666
$ws.Cells.Item(122, 10).Value = 2611.8333
$ws.Cells.Item(122, 12).Value = 7835.499899999999
$ws.Cells.Item(122, 14).Value = -12735.4999

# Row 125
$ws.Cells.Item(125, 8).Value = 969.2
$ws.Cells.Item(125, 9).Value = 895.1111
$ws.Cells.Item(125, 10).Value = 1636
$ws.Cells.Item(125, 11).Value = 8055.9999
$ws.Cells.Item(125, 12).Value = 14724
$ws.Cells.Item(125, 13).Value = -5595.9999
$ws.Cells.Item(125, 14).Value = -19644

# Row 129
$ws.Cells.Item(129, 8).Value = 1072.7646
$ws.Cells.Item(129, 10).Value = 1072.7646
$ws.Cells.Item(129, 12).Value = 3218.2938
$ws.Cells.Item(129, 14).Value = -13218.2938

# Row 134
$ws.Cells.Item(134, 8).Value = 42980
$ws.Cells.Item(134, 10).Value = 42980
$ws.Cells.Item(134, 12).Value = 42980
$ws.Cells.Item(134, 14).Value = -53120

# Row 137
$ws.Cells.Item(137, 8).Value = 1850.7838
$ws.Cells.Item(137, 9).Value = 1539
$ws.Cells.Item(137, 11).Value = 4617
$ws.Cells.Item(137, 13).Value = -2067

# Row 138
$ws.Cells.Item(138, 8).Value = 3600.1067
$ws.Cells.Item(138, 9).Value = 777.53845
$ws.Cells.Item(138, 10).Value = 5097.796
$ws.Cells.Item(138, 11).Value = 2332.61535
$ws.Cells.Item(138, 12).Value = 15293.388
$ws.Cells.Item(138, 13).Value = 2807.38465
$ws.Cells.Item(138, 14).Value = -25573.388

# Row 139
$ws.Cells.Item(139, 8).Value = 11952.167
$ws.Cells.Item(139, 10).Value = 11952.167
$ws.Cells.Item(139, 12).Value = 11952.167
$ws.Cells.Item(139, 14).Value = -22232.167

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 3937.524
$ws.Cells.Item(61, 9).Value = 1810.8
$ws.Cells.Item(61, 10).Value = 9254.333000000001
$ws.Cells.Item(61, 11).Value = 1810.8
$ws.Cells.Item(61, 12).Value = 9254.333000000001
$ws.Cells.Item(61, 13).Value = -1598.8
$ws.Cells.Item(61, 14).Value = -9678.333000000001

# Row 132
$ws.Cells.Item(132, 8).Value = 2759.2856
$ws.Cells.Item(132, 9).Value = 2208.0715
$ws.Cells.Item(132, 10).Value = 3861.7144
$ws.Cells.Item(132, 11).Value = 6624.2145
$ws.Cells.Item(132, 12).Value = 11585.1432
$ws.Cells.Item(132, 13).Value = -4094.2145
$ws.Cells.Item(132, 14).Value = -16645.1432

# Row 136
$ws.Cells.Item(136, 8).Value = 3937.524
$ws.Cells.Item(136, 9).Value = 1810.8
$ws.Cells.Item(136, 10).Value = 9254.333000000001
$ws.Cells.Item(136, 11).Value = 5432.4
$ws.Cells.Item(136, 12).Value = 27762.999
$ws.Cells.Item(136, 13).Value = -2882.4
$ws.Cells.Item(136, 14).Value = -32862.999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 2733.149
$ws.Cells.Item(31, 9).Value = 2070.0557
$ws.Cells.Item(31, 10).Value = 4903.273
$ws.Cells.Item(31, 11).Value = 2070.0557
$ws.Cells.Item(31, 12).Value = 4903.273
$ws.Cells.Item(31, 13).Value = -1775.0557
$ws.Cells.Item(31, 14).Value = -5493.273

# Row 34
$ws.Cells.Item(34, 8).Value = 2733.149
$ws.Cells.Item(34, 9).Value = 2070.0557
$ws.Cells.Item(34, 10).Value = 4903.273
$ws.Cells.Item(34, 11).Value = 2070.0557
$ws.Cells.Item(34, 12).Value = 4903.273
$ws.Cells.Item(34, 13).Value = -1868.0557
$ws.Cells.Item(34, 14).Value = -5307.273

# Row 132
$ws.Cells.Item(132, 8).Value = 3213.6667
$ws.Cells.Item(132, 9).Value = 2845.5
$ws.Cells.Item(132, 10).Value = 3950
$ws.Cells.Item(132, 11).Value = 8536.5
$ws.Cells.Item(132, 12).Value = 11850
$ws.Cells.Item(132, 13).Value = -6006.5
$ws.Cells.Item(132, 14).Value = -16910

# Row 134
$ws.Cells.Item(134, 8).Value = 2822.6667
$ws.Cells.Item(134, 9).Value = 1519.6
$ws.Cells.Item(134, 10).Value = 9338
$ws.Cells.Item(134, 11).Value = 4558.799999999999
$ws.Cells.Item(134, 12).Value = 28014
$ws.Cells.Item(134, 13).Value = -2023.799999999999
$ws.Cells.Item(134, 14).Value = -33084

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 714
$ws.Cells.Item(5, 9).Value = 588
$ws.Cells.Item(5, 10).Value = 912
$ws.Cells.Item(5, 11).Value = 1764
$ws.Cells.Item(5, 12).Value = 2736
$ws.Cells.Item(5, 13).Value = -1652
$ws.Cells.Item(5, 14).Value = -2960

# Row 131
$ws.Cells.Item(131, 8).Value = 4249.0347
$ws.Cells.Item(131, 9).Value = 531.6667
$ws.Cells.Item(131, 10).Value = 6873.0586
$ws.Cells.Item(131, 11).Value = 1595.0001
$ws.Cells.Item(131, 12).Value = 20619.1758
$ws.Cells.Item(131, 13).Value = 3444.9999
$ws.Cells.Item(131, 14).Value = -30699.1758

# Row 135
$ws.Cells.Item(135, 8).Value = 714
$ws.Cells.Item(135, 9).Value = 588
$ws.Cells.Item(135, 10).Value = 912
$ws.Cells.Item(135, 11).Value = 5292
$ws.Cells.Item(135, 12).Value = 8208
$ws.Cells.Item(135, 13).Value = -2757
$ws.Cells.Item(135, 14).Value = -13278

$ws = $wb.Worksheets.Item("GSM")
# Row 68
$ws.Cells.Item(68, 8).Value = 20000
$ws.Cells.Item(68, 10).Value = 20000
$ws.Cells.Item(68, 12).Value = 20000
$ws.Cells.Item(68, 14).Value = -21622

# Row 69
$ws.Cells.Item(69, 8).Value = 500010000
$ws.Cells.Item(69, 10).Value = 20000
$ws.Cells.Item(69, 12).Value = 20000
$ws.Cells.Item(69, 14).Value = -21498

# Row 71
$ws.Cells.Item(71, 8).Value = 20000
$ws.Cells.Item(71, 10).Value = 20000
$ws.Cells.Item(71, 12).Value = 60000
$ws.Cells.Item(71, 14).Value = -68112

# Row 72
$ws.Cells.Item(72, 8).Value = 500010000
$ws.Cells.Item(72, 10).Value = 20000
$ws.Cells.Item(72, 12).Value = 60000
$ws.Cells.Item(72, 14).Value = -67488

# Row 74
$ws.Cells.Item(74, 8).Value = 20087.334
$ws.Cells.Item(74, 10).Value = 20087.334
$ws.Cells.Item(74, 12).Value = 20087.334
$ws.Cells.Item(74, 14).Value = -21959.334

# Row 75
$ws.Cells.Item(75, 8).Value = 19748.5
$ws.Cells.Item(75, 10).Value = 19748.5
$ws.Cells.Item(75, 12).Value = 19748.5
$ws.Cells.Item(75, 14).Value = -21496.5

# Row 77
$ws.Cells.Item(77, 8).Value = 20087.334
$ws.Cells.Item(77, 10).Value = 20087.334
$ws.Cells.Item(77, 12).Value = 60262.00199999999
$ws.Cells.Item(77, 14).Value = -69622.00199999999

# Row 78
$ws.Cells.Item(78, 8).Value = 19748.5
$ws.Cells.Item(78, 10).Value = 19748.5
$ws.Cells.Item(78, 12).Value = 59245.5
$ws.Cells.Item(78, 14).Value = -67981.5

# Row 122
$ws.Cells.Item(122, 8).Value = 14287966
$ws.Cells.Item(122, 9).Value = 20001956
$ws.Cells.Item(122, 10).Value = 2990
$ws.Cells.Item(122, 11).Value = 60005868
$ws.Cells.Item(122, 12).Value = 8970
$ws.Cells.Item(122, 13).Value = -60003418
$ws.Cells.Item(122, 14).Value = -13870

$ws = $wb.Worksheets.Item("LTW")
# Row 30
$ws.Cells.Item(30, 8).Value = 1583
$ws.Cells.Item(30, 9).Value = 166
$ws.Cells.Item(30, 11).Value = 166
$ws.Cells.Item(30, 13).Value = -58

# Row 33
$ws.Cells.Item(33, 8).Value = 4000
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 10).Value = 4000
$ws.Cells.Item(33, 11).Value = 0
$ws.Cells.Item(33, 12).Value = 4000
$ws.Cells.Item(33, 13).Value = ""
$ws.Cells.Item(33, 14).Value = -4580

# Row 43
$ws.Cells.Item(43, 8).Value = 4750
$ws.Cells.Item(43, 10).Value = 4750
$ws.Cells.Item(43, 12).Value = 4750
$ws.Cells.Item(43, 14).Value = -5136

# Row 140
$ws.Cells.Item(140, 8).Value = 31807.666
$ws.Cells.Item(140, 10).Value = 31807.666
$ws.Cells.Item(140, 12).Value = 31807.666
$ws.Cells.Item(140, 14).Value = -42167.666

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Cells.Item(46, 8).Value = 50150
$ws.Cells.Item(46, 10).Value = 50150
$ws.Cells.Item(46, 12).Value = 50150
$ws.Cells.Item(46, 14).Value = -50612

# Row 107
$ws.Cells.Item(107, 8).Value = 547
$ws.Cells.Item(107, 9).Value = 523.6429000000001
$ws.Cells.Item(107, 10).Value = 601.5
$ws.Cells.Item(107, 11).Value = 1570.9287
$ws.Cells.Item(107, 12).Value = 1804.5
$ws.Cells.Item(107, 13).Value = 349.0712999999998
$ws.Cells.Item(107, 14).Value = -5644.5

# Row 134
$ws.Cells.Item(134, 8).Value = 50150
$ws.Cells.Item(134, 10).Value = 50150
$ws.Cells.Item(134, 12).Value = 150450
$ws.Cells.Item(134, 14).Value = -155520
